$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 2 into row 3 (keeps the same cell types/shared-string reuse
# that Excel would when copying a row), then overwrite the invoice-number
# cell with the new (buggy) external-invoice date string.
$ws.Range("A2:I2").Copy($ws.Range("A3:I3"))
$ws.Range("A3").Value = "1/03/2021asdasd"

# Column A widened (bestFit) to accommodate the longer text now in it.
$ws.Columns.Item(1).ColumnWidth = 15.6
